# Update the "kick" segment constraint table to be more permissive with
# kick varying between segments: change several "N" entries to "Y".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells (row 4-8, 10-12, 14 in column M; row 9 columns K and L; row 15 column K)
# that move from "N" to "Y".
$cellsToUpdate = @("M4", "M5", "M6", "M7", "M8", "K9", "L9", "M10", "M11", "M12", "M14", "K15")

foreach ($cell in $cellsToUpdate) {
    $ws.Range($cell).Value = "Y"
}
